# DB_Setup.xlsx edit script
# Commit message: "need to insert facs data and after implement unique forms for users"
#
# Semantic changes applied:
#  1. Informatii_Facultati column list (G7:H16): a new "ID_fac - FK" row is
#     inserted at the top, pushing Durata licenta / Ultima medie admitere /
#     Tip admitere down by one row; the large blank block below is unmerged
#     (only the new top row stays merged).
#  2. Domenii de studiu column list (J4:K9): "Nume" moves up from J8 to J7
#     (now merged+centered there); the table loses its last (blank) row
#     J10:K10.
#  3. Users table (rows 16-22 -> 16-19): replaced by a new, simpler 3-row
#     column list (username - PK / tip / passwd); the old
#     Username/Password/Tip-user layout and the separate tip-user lookup
#     table (admin/elev/student/profesor/consilier cariera) are removed.
#  4. Selection/view state updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Users table (rows 16-22 -> 16-19)
# ---------------------------------------------------------------------
$ws.Range("C17").Value = "username - PK"
$ws.Range("C18").Value = "tip "
$ws.Range("C19").Value = "passwd"
$ws.Range("C20").Value = ""

$ws.Range("E16").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("F21").Value = ""
$ws.Range("F22").Value = ""

# ---------------------------------------------------------------------
# 2. Domenii de studiu columns (J4:K9)
# ---------------------------------------------------------------------
# Drop the trailing blank row of the table.
$ws.Range("J10:K10").Clear() | Out-Null
$ws.Range("J10:K10").UnMerge() | Out-Null

# "Nume" moves from J8 up to J7, which becomes its own merged cell; J7 is
# already centre-aligned from its old (empty) style, so no extra formatting
# call is needed. J8 becomes blank (matching the row it vacated).
$ws.Range("J8").Value = ""
$ws.Range("J7:K7").Merge() | Out-Null
$ws.Range("J7").Value = "Nume"

# ---------------------------------------------------------------------
# 3. Informatii_Facultati columns (G7:H16)
# ---------------------------------------------------------------------
# Break up the old per-row merges so we can restructure the block.
$ws.Range("G8:H8").UnMerge() | Out-Null
$ws.Range("G9:H9").UnMerge() | Out-Null
$ws.Range("G10:H10").UnMerge() | Out-Null
$ws.Range("G11:H16").UnMerge() | Out-Null

# Shift the existing three entries down one row and insert the new one on top.
$ws.Range("G7").Value = "ID_fac - FK"
$ws.Range("G8").Value = "Durata licenta"
$ws.Range("G9").Value = "Ultima medie admitere"
$ws.Range("G10").Value = "Tip admitere"

# G7:H7 remains its own merged/centered header-style row (style unchanged).
# The remaining rows (G8:H16) become individual, unmerged, generically
# aligned cells.
$ws.Range("G8:H16").HorizontalAlignment = 1   # xlGeneral

# ---------------------------------------------------------------------
# 4. View/selection state
# ---------------------------------------------------------------------
$ws.Range("E13").Select() | Out-Null
